$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A2" = 0.8863028663560529
    "B2" = 0.8866365902508907
    "C2" = 0.886970215846263
    "D2" = 0.8873037432159071
    "E2" = 0.8873037432159071
    "F2" = 0.8873037432159071
    "G2" = 0.8873037432159071
    "H2" = 0.8867366882475899
    "I2" = 0.8873037432159071
    "J2" = 0.8889699092528377
    "K2" = 0.885635123372578
    "L2" = 0.8785663091447938
    "M2" = 0.8714529261598966
    "N2" = 0.8642942513898125
    "O2" = 0.8867366882475899
    "P2" = 0.8867366882475899
    "Q2" = 0.8811845799953736
    "R2" = 0.8845324856007161
    "S2" = 0.8878705145246762
    "T2" = 0.8911987406678564
    "U2" = 0.8945172370459996
    "V2" = 0.9087757896185328
    "W2" = 0.905498601466514
    "X2" = 0.9087757896185328
    "Y2" = 0.8911987406678564
    "Z2" = 0.9007518923840988
    "AA2" = 0.9102256793195173
    "AB2" = 0.9196217600848949
    "AC2" = 0.9196217600848949
    "AD2" = 0.9163754624134822
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
